# Edit script: shifts the weekly Betarraga price records down by 3 rows
# (rows 441-556 -> 444-559), inserts a brand-new date block (44508) at
# rows 441-443, and extends the sheet dimension to A1:R559.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 441
$lastRow  = 556

# Columns that carry the per-record data that shifts with the rows.
# (D = Fecha, I = Calidad, J = Volumen, K = Precio minimo, L = Precio maximo,
#  M = Precio promedio ponderado, P = Precio $/Kg)

# 1) Snapshot every "before" value for the block we are about to shift into
#    plain arrays first, so later writes can never clobber data we still
#    need to read.
$D = @()
$I = @()
$J = @()
$K = @()
$L = @()
$M = @()
$P = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $D += $ws.Cells.Item($r, 4).Value2()
    $I += $ws.Cells.Item($r, 9).Value2()
    $J += $ws.Cells.Item($r, 10).Value2()
    $K += $ws.Cells.Item($r, 11).Value2()
    $L += $ws.Cells.Item($r, 12).Value2()
    $M += $ws.Cells.Item($r, 13).Value2()
    $P += $ws.Cells.Item($r, 16).Value2()
}

# 2) The last 3 rows of the old block (554-556) become 3 brand new rows
#    (557-559). Build those new rows explicitly from the static column
#    values (identical for every record in this block) plus the shifted
#    data snapshotted above, rather than Copy/PasteSpecial (which would
#    otherwise register a brand-new duplicate number-format style).
$dateFormat = $ws.Cells.Item($firstRow, 4).NumberFormat()

$staticA = $ws.Cells.Item($firstRow, 1).Value2()
$staticB = $ws.Cells.Item($firstRow, 2).Value2()
$staticC = $ws.Cells.Item($firstRow, 3).Value2()
$staticE = $ws.Cells.Item($firstRow, 5).Value2()
$staticF = $ws.Cells.Item($firstRow, 6).Value2()
$staticG = $ws.Cells.Item($firstRow, 7).Value2()
$staticH = $ws.Cells.Item($firstRow, 8).Value2()
$staticN = $ws.Cells.Item($firstRow, 14).Value2()
$staticO = $ws.Cells.Item($firstRow, 15).Value2()
$staticQ = $ws.Cells.Item($firstRow, 17).Value2()
$staticR = $ws.Cells.Item($firstRow, 18).Value2()

for ($k = 0; $k -lt 3; $k++) {
    $srcIdx  = $lastRow - 2 + $k - $firstRow   # snapshot index for rows 554, 555, 556
    $dstRow  = $lastRow + 1 + $k                # 557, 558, 559

    $ws.Cells.Item($dstRow, 1).Value2  = $staticA
    $ws.Cells.Item($dstRow, 2).Value2  = $staticB
    $ws.Cells.Item($dstRow, 3).Value2  = $staticC
    $ws.Cells.Item($dstRow, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($dstRow, 4).Value2  = $D[$srcIdx]
    $ws.Cells.Item($dstRow, 5).Value2  = $staticE
    $ws.Cells.Item($dstRow, 6).Value2  = $staticF
    $ws.Cells.Item($dstRow, 7).Value2  = $staticG
    $ws.Cells.Item($dstRow, 8).Value2  = $staticH
    $ws.Cells.Item($dstRow, 9).Value2  = $I[$srcIdx]
    $ws.Cells.Item($dstRow, 10).Value2 = $J[$srcIdx]
    $ws.Cells.Item($dstRow, 11).Value2 = $K[$srcIdx]
    $ws.Cells.Item($dstRow, 12).Value2 = $L[$srcIdx]
    $ws.Cells.Item($dstRow, 13).Value2 = $M[$srcIdx]
    $ws.Cells.Item($dstRow, 14).Value2 = $staticN
    $ws.Cells.Item($dstRow, 15).Value2 = $staticO
    $ws.Cells.Item($dstRow, 16).Value2 = $P[$srcIdx]
    $ws.Cells.Item($dstRow, 17).Value2 = $staticQ
    $ws.Cells.Item($dstRow, 18).Value2 = $staticR
}

# 3) Shift rows 444-556 <- rows 441-553 (i.e. new[r] = old[r-3]) using the
#    snapshot arrays captured in step 1 (order no longer matters).
for ($r = ($firstRow + 3); $r -le $lastRow; $r++) {
    $idx = $r - 3 - $firstRow
    $ws.Cells.Item($r, 4).Value2  = $D[$idx]
    $ws.Cells.Item($r, 9).Value2  = $I[$idx]
    $ws.Cells.Item($r, 10).Value2 = $J[$idx]
    $ws.Cells.Item($r, 11).Value2 = $K[$idx]
    $ws.Cells.Item($r, 12).Value2 = $L[$idx]
    $ws.Cells.Item($r, 13).Value2 = $M[$idx]
    $ws.Cells.Item($r, 16).Value2 = $P[$idx]
}

# 4) Rows 441-443 get the brand new "44508" date block (Primera/Segunda/Tercera).
$ws.Cells.Item(441, 4).Value2  = 44508
$ws.Cells.Item(441, 9).Value2  = "Primera"
$ws.Cells.Item(441, 10).Value2 = 60000
$ws.Cells.Item(441, 11).Value2 = 90
$ws.Cells.Item(441, 12).Value2 = 100
$ws.Cells.Item(441, 13).Value2 = 94
$ws.Cells.Item(441, 16).Value2 = 94

$ws.Cells.Item(442, 4).Value2  = 44508
$ws.Cells.Item(442, 9).Value2  = "Segunda"
$ws.Cells.Item(442, 10).Value2 = 62000
$ws.Cells.Item(442, 11).Value2 = 75
$ws.Cells.Item(442, 12).Value2 = 80
$ws.Cells.Item(442, 13).Value2 = 78
$ws.Cells.Item(442, 16).Value2 = 78

$ws.Cells.Item(443, 4).Value2  = 44508
$ws.Cells.Item(443, 9).Value2  = "Tercera"
$ws.Cells.Item(443, 10).Value2 = 9500
$ws.Cells.Item(443, 11).Value2 = 60
$ws.Cells.Item(443, 12).Value2 = 60
$ws.Cells.Item(443, 13).Value2 = 60
$ws.Cells.Item(443, 16).Value2 = 60
